$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns I ("I0") and J ("IF") ---
# Copy formatting (bold font, border, centered alignment) from the existing
# header cell H1 so I1/J1 match the style used by the other header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-40: values for columns I and J ---
$I = @(6,5,6,6,7,7,7,8,7,8,10,5,8,7,8,7,7,7,7,7,6,7,10,7,10,7,7,7,8,7,7,9,6,8,6,8,7,7,8)
$J = @(6,5,6,6,7,7,7,8,7,8,10,6,8,7,8,7,7,7,7,7,6,7,10,7,10,7,8,7,8,7,7,9,7,9,6,8,7,7,8)

for ($n = 0; $n -lt $I.Length; $n++) {
    $row = $n + 2
    $ws.Cells.Item($row, 9).Value = $I[$n]
    $ws.Cells.Item($row, 10).Value = $J[$n]
}
